$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 11

$ws.Cells.Item($row, 1).Value = 45945.83580382937
$ws.Cells.Item($row, 2).Value = "swency"
$ws.Cells.Item($row, 3).Value = "There's a cylinder burst near my area, at nehru nagar, coimbatore. Immediate Help is required"
$ws.Cells.Item($row, 4).Value = "Water"
$ws.Cells.Item($row, 5).Value = "High"
$ws.Cells.Item($row, 6).Value = "area, at, at nehru, near my"
$ws.Cells.Item($row, 7).Value = "COMP_20251015_200333_5845"
$ws.Cells.Item($row, 8).Value = "complaint_COMP_20251015_200333_5845_20251015_200333_20251015_200333.jpeg"

$ws.Cells.Item($row, 1).NumberFormat = $ws.Cells.Item($row - 1, 1).NumberFormat
